$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Natmi LR-pairs output for Cd84-Cd84: adds the "sCs" cluster alongside the
# existing "ECs" cluster and refreshes the ECs-ECs edge stats (row 2), producing
# all four ECs/sCs sender-target combinations (rows 2-5). Columns A-T are:
# Sending cluster, Ligand symbol, Receptor symbol, Target cluster,
# Ligand-expressing cells, Ligand detection rate, Ligand average expression value,
# Ligand total expression value, Ligand derived specificity (avg), Ligand derived specificity (total),
# Receptor-expressing cells, Receptor detection rate, Receptor average expression value,
# Receptor total expression value, Receptor derived specificity (avg), Receptor derived specificity (total),
# Edge average expression weight, Edge total expression weight,
# Edge average expression derived specificity, Edge total expression derived specificity
$dataRows = @(
    @("ECs","Cd84","Cd84","ECs",3,1,56.376102,169.128306,0.9999051815767483,0.9999051815767483,3,1,56.376102,169.128306,0.9999051815767483,0.9999051815767483,3178.264876714404,28604.38389042964,0.9998103721440301,0.9998103721440301),
    @("ECs","Cd84","Cd84","sCs",3,1,56.376102,169.128306,0.9999051815767483,0.9999051815767483,1,0.3333333333333333,0.005346,0.016038,0.00009481842325156316,0.00009481842325156316,0.301386641292,2.712479771628,0.00009480943271817525,0.00009480943271817525),
    @("sCs","Cd84","Cd84","ECs",1,0.3333333333333333,0.005346,0.016038,0.00009481842325156316,0.00009481842325156316,3,1,56.376102,169.128306,0.9999051815767483,0.9999051815767483,0.301386641292,2.712479771628,0.00009480943271817525,0.00009480943271817525),
    @("sCs","Cd84","Cd84","sCs",1,0.3333333333333333,0.005346,0.016038,0.00009481842325156316,0.00009481842325156316,1,0.3333333333333333,0.005346,0.016038,0.00009481842325156316,0.00009481842325156316,0.000028579716,0.000257217444,0.000000008990533387912575,0.000000008990533387912575),
)

$startRow = 2
for ($i = 0; $i -lt $dataRows.Length; $i++) {
    $rowNum = $startRow + $i
    $rowValues = $dataRows[$i]
    for ($col = 1; $col -le $rowValues.Length; $col++) {
        $ws.Cells.Item($rowNum, $col).Value2 = $rowValues[$col - 1]
    }
}
